# Adding the changes we made on may 9th
# The accelerometer sample window rolled forward: the oldest 5 rows of
# readings were dropped and 10 new rows of x/y/z samples were appended,
# so the sheet now holds 30 data rows (A2:C31) instead of 20 (A2:C21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = New-Object "object[,]" 30,3
$newValues[0,0] = 4.936917678169539
$newValues[0,1] = -12.71771009072015
$newValues[0,2] = 3.41959030731864
$newValues[1,0] = 11.45307619675343
$newValues[1,1] = -9.954130877619217
$newValues[1,2] = 7.623632503592415
$newValues[2,0] = -3.206811956737363
$newValues[2,1] = -14.22782378611358
$newValues[2,2] = -3.585576596467396
$newValues[3,0] = -12.97183887854877
$newValues[3,1] = -19.3549799711808
$newValues[3,2] = -15.49820016778038
$newValues[4,0] = 6.874614751857277
$newValues[4,1] = -17.72338112540867
$newValues[4,2] = -7.79767358821371
$newValues[5,0] = 27.40819798345139
$newValues[5,1] = -14.00006047539093
$newValues[5,2] = 1.457801404206576
$newValues[6,0] = 40.11718940734848
$newValues[6,1] = -4.808408617973335
$newValues[6,2] = 16.45111835002897
$newValues[7,0] = -6.829358805780815
$newValues[7,1] = -18.71223431048175
$newValues[7,2] = 7.999558770138284
$newValues[8,0] = -8.796926493230012
$newValues[8,1] = -58.53746407446633
$newValues[8,2] = 6.315518513969705
$newValues[9,0] = 3.485328860904898
$newValues[9,1] = -5.610776631728454
$newValues[9,2] = -5.665500184764032
$newValues[10,0] = -16.32176755822241
$newValues[10,1] = -11.76148359671859
$newValues[10,2] = -5.036116931749546
$newValues[11,0] = 16.76886541947075
$newValues[11,1] = -55.28290149439937
$newValues[11,2] = 16.80483585855235
$newValues[12,0] = -17.68672215420252
$newValues[12,1] = -8.062301013780889
$newValues[12,2] = 3.262691987597485
$newValues[13,0] = 6.037598153819246
$newValues[13,1] = -26.96122758284862
$newValues[13,2] = 22.58520386530009
$newValues[14,0] = -27.64734548071142
$newValues[14,1] = -30.63590854147206
$newValues[14,2] = -13.72450681354674
$newValues[15,0] = -0.3365890254143977
$newValues[15,1] = -9.820995450019772
$newValues[15,2] = -11.5744883495828
$newValues[16,0] = -9.280692992003011
$newValues[16,1] = -14.45616371476128
$newValues[16,2] = -11.51578338249869
$newValues[17,0] = -4.598872557930385
$newValues[17,1] = -17.49396556356679
$newValues[17,2] = 6.897694048674124
$newValues[18,0] = -11.20024363890936
$newValues[18,1] = 1.505196239637264
$newValues[18,2] = 13.08335323955702
$newValues[19,0] = -36.36955037324309
$newValues[19,1] = -40.35419501428967
$newValues[19,2] = 41.72007127430098
$newValues[20,0] = -64.13763353098994
$newValues[20,1] = -54.61329487095708
$newValues[20,2] = 20.88871420984682
$newValues[21,0] = -48.54061980869488
$newValues[21,1] = -15.36169694817591
$newValues[21,2] = -18.69812476116687
$newValues[22,0] = -9.763155895730725
$newValues[22,1] = -17.02286973725198
$newValues[22,2] = -7.877094351727034
$newValues[23,0] = -12.50246284319005
$newValues[23,1] = -25.3144741265671
$newValues[23,2] = 2.428819822228345
$newValues[24,0] = -12.77425661294355
$newValues[24,1] = -3.042295404102455
$newValues[24,2] = 24.03419656857195
$newValues[25,0] = -11.1896470111351
$newValues[25,1] = -11.74091952780091
$newValues[25,2] = 15.447055526402
$newValues[26,0] = -33.40386452882192
$newValues[26,1] = -74.93744767230478
$newValues[26,2] = 35.67094943834376
$newValues[27,0] = -55.71352932764142
$newValues[27,1] = -17.14117759207022
$newValues[27,2] = -22.77800974638568
$newValues[28,0] = -7.357292139011771
$newValues[28,1] = -6.564726891724789
$newValues[28,2] = -1.292621791362762
$newValues[29,0] = -26.65620994567871
$newValues[29,1] = -36.59538269042969
$newValues[29,2] = -1.42856240272522

$target = $ws.Range("A2").Resize($newValues.GetLength(0), $newValues.GetLength(1))
$target.Value = $newValues
